$d = $word.ActiveDocument

$replacements = @(
    @("817×9=", "988×5="),
    @("383×3=", "824×2="),
    @("145×7=", "623×9="),
    @("563×6=", "567×6="),
    @("157×8=", "885×8="),
    @("482×4=", "839×3="),
    @("404×5=", "696×4="),
    @("767×4=", "366×3="),
    @("601×3=", "671×8="),
    @("561×9=", "766×8="),
    @("856×3=", "418×2="),
    @("446×3=", "761×7="),
    @("662×7=", "263×7="),
    @("649×5=", "619×7="),
    @("348×3=", "526×7="),
    @("575×3=", "850×7="),
    @("687×3=", "247×9="),
    @("720×5=", "450×9="),
    @("122×8=", "869×9="),
    @("695×3=", "873×4="),
    @("990×5=", "908×7="),
    @("420×5=", "194×7="),
    @("775×9=", "194×2="),
    @("480×5=", "122×4="),
    @("918×2=", "743×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
